# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 89 (pushing the
# existing rows 89-117 down to 90-118) for
# "Feria Lagunitas de Puerto Montt - Haba".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 89, shifting rows 89:117 down to 90:118.
$ws.Rows.Item(89).Insert()

# Populate the newly inserted row 89 with the new weekly record.
$ws.Cells.Item(89, 1).Value  = 4
$ws.Cells.Item(89, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(89, 3).Value  = "Los Lagos"
$ws.Cells.Item(89, 4).Value  = 44876
$ws.Cells.Item(89, 5).Value  = 10
$ws.Cells.Item(89, 6).Value  = 100112026
$ws.Cells.Item(89, 7).Value  = "Haba"
$ws.Cells.Item(89, 8).Value  = "Sin especificar"
$ws.Cells.Item(89, 9).Value  = "Primera"
$ws.Cells.Item(89, 10).Value = 160
$ws.Cells.Item(89, 11).Value = 12000
$ws.Cells.Item(89, 12).Value = 13000
$ws.Cells.Item(89, 13).Value = 12500
$ws.Cells.Item(89, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(89, 15).Value = "Región del Maule"
$ws.Cells.Item(89, 16).Value = 500
$ws.Cells.Item(89, 17).Value = 25
$ws.Cells.Item(89, 18).Value = "Hortaliza"
